$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, pushing the existing data (old rows 69-107)
# down to rows 70-108.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with a new weekly observation
# (same market/category metadata as the row that used to be there, but a
# new date and recomputed prices).
$ws.Cells.Item(69, 1).Value = 6
$ws.Cells.Item(69, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(69, 3).Value = "Metropolitana"
$ws.Cells.Item(69, 4).Value = 44460
$ws.Cells.Item(69, 5).Value = 13
$ws.Cells.Item(69, 6).Value = 100112029
$ws.Cells.Item(69, 7).Value = "Orégano"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 29
$ws.Cells.Item(69, 11).Value = 9000
$ws.Cells.Item(69, 12).Value = 10000
$ws.Cells.Item(69, 13).Value = 9483
$ws.Cells.Item(69, 14).Value = "$/docena de atados"
$ws.Cells.Item(69, 15).Value = "Región Metropolitana"
$ws.Cells.Item(69, 16).Value = 3161
$ws.Cells.Item(69, 17).Value = 3
$ws.Cells.Item(69, 18).Value = "Hortaliza"
